$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at row 270 (pushes old rows 270..301 down to 271..302) ---
$ws.Rows.Item(270).Insert()

$ws.Range("A270").Value = 10
$ws.Range("B270").Value = "Vega Modelo de Temuco"
$ws.Range("C270").Value = "La Araucanía"
$ws.Range("D270").Value = 44748
$ws.Range("E270").Value = 9
$ws.Range("F270").Value = 100112001
$ws.Range("G270").Value = "Berenjena"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 45
$ws.Range("K270").Value = 12000
$ws.Range("L270").Value = 13000
$ws.Range("M270").Value = 12444
$ws.Range("N270").Value = "$/caja 60 unidades"
$ws.Range("O270").Value = "Región de Arica y Parinacota"
$ws.Range("P270").Value = 207
$ws.Range("Q270").Value = 60
$ws.Range("R270").Value = "Hortaliza"

# --- Insert second new row at row 292 (pushes the (already shifted) old rows 291..302 down to 292..303) ---
$ws.Rows.Item(292).Insert()

$ws.Range("A292").Value = 10
$ws.Range("B292").Value = "Vega Modelo de Temuco"
$ws.Range("C292").Value = "La Araucanía"
$ws.Range("D292").Value = 44747
$ws.Range("E292").Value = 9
$ws.Range("F292").Value = 100112001
$ws.Range("G292").Value = "Berenjena"
$ws.Range("H292").Value = "Sin especificar"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 45
$ws.Range("K292").Value = 12000
$ws.Range("L292").Value = 12000
$ws.Range("M292").Value = 12000
$ws.Range("N292").Value = "$/caja 60 unidades"
$ws.Range("O292").Value = "Región de Arica y Parinacota"
$ws.Range("P292").Value = 200
$ws.Range("Q292").Value = 60
$ws.Range("R292").Value = "Hortaliza"
